$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 232.9333528440843
$ws.Range("R2").Value = 2096.400175596759
$ws.Range("S2").Value = 0.004898908148186765
$ws.Range("T2").Value = 0.004898908148186765

# Row 3
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 15787.49651441522
$ws.Range("R3").Value = 142087.468629737
$ws.Range("S3").Value = 0.3320327225346212
$ws.Range("T3").Value = 0.3320327225346212

# Row 4
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 5049.635733903989
$ws.Range("R4").Value = 45446.72160513589
$ws.Range("S4").Value = 0.1062007709078728
$ws.Range("T4").Value = 0.1062007709078728

# Row 5
$ws.Range("G5").Value = 169.915657
$ws.Range("H5").Value = 509.746971
$ws.Range("I5").Value = 0.4441184931734509
$ws.Range("J5").Value = 0.4441184931734509
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 46.886696308237
$ws.Range("R5").Value = 421.980266774133
$ws.Range("S5").Value = 0.000986091582770148
$ws.Range("T5").Value = 0.000986091582770148

# Row 6
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 93.74359854894901
$ws.Range("R6").Value = 843.692386940541
$ws.Range("S6").Value = 0.001971556555403175
$ws.Range("T6").Value = 0.001971556555403175

# Row 7
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.1336259572377962
$ws.Range("T7").Value = 0.1336259572377962

# Row 8
$ws.Range("I8").Value = 0.1787346690539575
$ws.Range("J8").Value = 0.1787346690539575
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 2032.216594479612
$ws.Range("R8").Value = 18289.94935031651
$ws.Range("S8").Value = 0.04274030452066867
$ws.Range("T8").Value = 0.04274030452066867

# Row 9
$ws.Range("I9").Value = 0.1787346690539575
$ws.Range("J9").Value = 0.1787346690539575
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 18.869464515663
$ws.Range("R9").Value = 169.825180640967
$ws.Range("S9").Value = 0.0003968507400894956
$ws.Range("T9").Value = 0.0003968507400894955

# Row 10
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 73.03293491759533
$ws.Range("R10").Value = 657.2964142583579
$ws.Range("S10").Value = 0.001535982870573649
$ws.Range("T10").Value = 0.001535982870573649

# Row 11
$ws.Range("G11").Value = 53.27463399999999
$ws.Range("H11").Value = 159.823902
$ws.Range("I11").Value = 0.1392470275793777
$ws.Range("J11").Value = 0.1392470275793778
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 4949.944657435227
$ws.Range("R11").Value = 44549.50191691704
$ws.Range("S11").Value = 0.1041041307279616
$ws.Range("T11").Value = 0.1041041307279617

# Row 12
$ws.Range("G12").Value = 53.27463399999999
$ws.Range("H12").Value = 159.823902
$ws.Range("I12").Value = 0.1392470275793777
$ws.Range("J12").Value = 0.1392470275793778
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 1583.241358134856
$ws.Range("R12").Value = 14249.1722232137
$ws.Range("S12").Value = 0.03329773901080092
$ws.Range("T12").Value = 0.03329773901080092

# Row 13
$ws.Range("G13").Value = 53.27463399999999
$ws.Range("H13").Value = 159.823902
$ws.Range("I13").Value = 0.1392470275793777
$ws.Range("J13").Value = 0.1392470275793778
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 14.700655780594
$ws.Range("R13").Value = 132.305902025346
$ws.Range("S13").Value = 0.0003091749700415209
$ws.Range("T13").Value = 0.0003091749700415209

# Row 14
$ws.Range("G14").Value = 91.01828266666666
$ws.Range("H14").Value = 273.054848
$ws.Range("I14").Value = 0.2378998101932138
$ws.Range("J14").Value = 0.2378998101932138
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 124.7748096083769
$ws.Range("R14").Value = 1122.973286475392
$ws.Range("S14").Value = 0.002624185519229105
$ws.Range("T14").Value = 0.002624185519229105

# Row 15
$ws.Range("G15").Value = 91.01828266666666
$ws.Range("H15").Value = 273.054848
$ws.Range("I15").Value = 0.2378998101932138
$ws.Range("J15").Value = 0.2378998101932138
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 8456.847624984079
$ws.Range("R15").Value = 76111.6286248567
$ws.Range("S15").Value = 0.1778591139146115
$ws.Range("T15").Value = 0.1778591139146115

# Row 16
$ws.Range("G16").Value = 91.01828266666666
$ws.Range("H16").Value = 273.054848
$ws.Range("I16").Value = 0.2378998101932138
$ws.Range("J16").Value = 0.2378998101932138
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 2704.925377136811
$ws.Range("R16").Value = 24344.3283942313
$ws.Range("S16").Value = 0.05688829361917292
$ws.Range("T16").Value = 0.05688829361917291

# Row 17
$ws.Range("G17").Value = 91.01828266666666
$ws.Range("H17").Value = 273.054848
$ws.Range("I17").Value = 0.2378998101932138
$ws.Range("J17").Value = 0.2378998101932138
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 25.11567593732267
$ws.Range("R17").Value = 226.041083435904
$ws.Range("S17").Value = 0.0003091749700415209
$ws.Range("T17").Value = 0.0003091749700415209

Write-Host "Updated cells for rows 2-17"